# Append the latest EUR->ARS quote (2025-10-09T21:22:53Z) as a new row
# at the bottom of the rate-history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 68

# Column A holds a date-like string ("2025-10-09"). Excel would normally
# auto-convert a bare "2025-10-09" into a date serial number. Prefixing
# with a leading apostrophe forces it to be entered as literal text (as
# all the other rows in this column already are), and ClearFormats()
# afterwards drops the transient "quote prefix" cell format so the cell
# is left with the same default styling as its neighbours.
$ws.Range("A" + $newRow).Value = "'2025-10-09"
$ws.Range("A" + $newRow).ClearFormats()

# Column B ("21:22:53") and column C (the rate string) are plain text
# that Excel's type-inference leaves alone, so they can be set directly.
$ws.Range("B" + $newRow).Value = "21:22:53"
$ws.Range("C" + $newRow).Value = "1.00 EUR = 1,758.1419"
